# Apply cryptos.xlsx data refresh (GitHub Actions update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "68.402.56"
$ws.Cells.Item(2, 5).Value = "  +0.06%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.900.04"
$ws.Cells.Item(3, 5).Value = "  +0.13%  "

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "602.12"
$ws.Cells.Item(5, 5).Value = "  +0.04%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "169.06"
$ws.Cells.Item(6, 5).Value = "  +1.37%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "3.897.20"
$ws.Cells.Item(7, 5).Value = "  +0.05%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.02%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.531"
$ws.Cells.Item(9, 5).Value = "  +0.63%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -0.61%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "6.45"
$ws.Cells.Item(11, 5).Value = "  +0.28%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.459"
$ws.Cells.Item(12, 5).Value = "  -0.20%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000255"
$ws.Cells.Item(13, 5).Value = "  +0.67%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "37.23"
$ws.Cells.Item(14, 5).Value = "  -0.23%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "4.554.75"
$ws.Cells.Item(15, 5).Value = "  +0.17%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.900.31"
$ws.Cells.Item(16, 5).Value = "  +0.02%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "68.454.71"
$ws.Cells.Item(17, 5).Value = "  -0.01%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "18.26"
$ws.Cells.Item(18, 5).Value = "  +5.62%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.40"
$ws.Cells.Item(19, 5).Value = "  -0.87%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.10%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  -2.07%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "471.37"
$ws.Cells.Item(22, 5).Value = "  -3.63%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +2.09%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.0000166"
$ws.Cells.Item(24, 5).Value = "  +0.26%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "83.73"
$ws.Cells.Item(25, 5).Value = "  -1.09%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.83%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +1.69%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.08%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "10.00"
$ws.Cells.Item(29, 5).Value = "  -1.48%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +1.43%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.050.63"
$ws.Cells.Item(31, 5).Value = "  +0.17%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +2.37%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -2.22%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "31.47"
$ws.Cells.Item(34, 5).Value = "  -0.92%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "9.49"
$ws.Cells.Item(35, 5).Value = "  +1.94%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.868.22"
$ws.Cells.Item(36, 5).Value = "  +0.56%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.106"
$ws.Cells.Item(37, 5).Value = "  -1.03%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.66"
$ws.Cells.Item(38, 5).Value = "  +16.10%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.141"
$ws.Cells.Item(39, 5).Value = "  +1.72%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Mantle"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.03"
$ws.Cells.Item(40, 5).Value = "  -0.50%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.94"
$ws.Cells.Item(41, 5).Value = "  +0.29%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.999"
$ws.Cells.Item(42, 5).Value = "  -0.08%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.68%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.000302"
$ws.Cells.Item(44, 5).Value = "  +12.50%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "426.95"
$ws.Cells.Item(45, 5).Value = "  -1.13%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +0.00%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +1.41%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "47.17"
$ws.Cells.Item(49, 5).Value = "  -2.00%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "27.11"
$ws.Cells.Item(50, 5).Value = "  +6.32%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "143.44"
$ws.Cells.Item(51, 5).Value = "  +0.49%  "

